$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "83÷4="  # was 16÷8=
$t.Cell(1, 2).Range.Text = "85÷5="  # was 68÷8=
$t.Cell(1, 3).Range.Text = "46÷2="  # was 21÷7=
$t.Cell(1, 4).Range.Text = "60÷5="  # was 34÷9=
$t.Cell(1, 5).Range.Text = "47÷4="  # was 50÷7=

$t.Cell(5, 1).Range.Text = "82÷2="  # was 84÷4=
$t.Cell(5, 2).Range.Text = "72÷5="  # was 10÷2=
$t.Cell(5, 3).Range.Text = "46÷8="  # was 74÷2=
$t.Cell(5, 4).Range.Text = "48÷4="  # was 12÷6=
$t.Cell(5, 5).Range.Text = "36÷3="  # was 38÷9=

$t.Cell(9, 1).Range.Text = "90÷8="  # was 85÷6=
$t.Cell(9, 2).Range.Text = "50÷3="  # was 45÷7=
$t.Cell(9, 3).Range.Text = "36÷4="  # was 64÷8=
$t.Cell(9, 4).Range.Text = "48÷3="  # was 25÷5=
$t.Cell(9, 5).Range.Text = "42÷9="  # was 44÷3=

$t.Cell(13, 1).Range.Text = "31÷3="  # was 86÷8=
$t.Cell(13, 2).Range.Text = "41÷4="  # was 37÷9=
$t.Cell(13, 3).Range.Text = "16÷5="  # was 29÷3=
$t.Cell(13, 4).Range.Text = "78÷8="  # was 61÷2=
$t.Cell(13, 5).Range.Text = "83÷5="  # was 56÷3=

$t.Cell(17, 1).Range.Text = "82÷8="  # was 76÷9=
$t.Cell(17, 2).Range.Text = "31÷3="  # was 65÷9=
$t.Cell(17, 3).Range.Text = "36÷3="  # was 56÷3=
$t.Cell(17, 4).Range.Text = "58÷7="  # was 44÷6=
$t.Cell(17, 5).Range.Text = "55÷4="  # was 19÷5=
